# Insert a new data row right after the current row 44 (i.e. it becomes the
# new row 45), pushing every existing row from 45..125 down by one
# (45->46, ..., 125->126). Then populate the newly-inserted row 45 with the
# new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 45..125 down to 46..126 and leave a blank row 45 behind.
$ws.Rows(45).Insert()

# Fill the new row 45 with the new "Crespo record" observation.
$ws.Range("A45").Value = 7
$ws.Range("B45").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C45").Value = "Ñuble"
$ws.Range("D45").Value = 44469
$ws.Range("E45").Value = 16
$ws.Range("F45").Value = 100112006
$ws.Range("G45").Value = "Repollo"
$ws.Range("H45").Value = "Crespo record"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 300
$ws.Range("K45").Value = 600
$ws.Range("L45").Value = 650
$ws.Range("M45").Value = 625
$ws.Range("N45").Value = "`$/unidad"
$ws.Range("O45").Value = "Provincia de Diguillín"
$ws.Range("P45").Value = 625
$ws.Range("Q45").Value = 1
$ws.Range("R45").Value = "Hortaliza"

# Make sure the date cell keeps the same date-number-format style used by
# every other row in column D.
$ws.Range("D45").NumberFormat = $ws.Range("D46").NumberFormat
